{"js": "// Remove the trailing \"Ver no Jupiter / Salvar em pdf / Salvar em docx\" line,\n// the \"\u00a9 2020 ... Creative Commons Attribution\" footer line, and the blank\n// paragraph that separated them from the bibliography, while leaving the\n// rest of the document (including the FLEMMING bibliography entry and the\n// paragraph(s) that follow the footer) untouched.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the bibliography entry that immediately precedes the block we\n// need to drop.\nlet anchor = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"FLEMMING, Diva M.; GON\u00c7ALVES, Mirian B.\") !== -1) {\n    anchor = p;\n    break;\n  }\n}\n\nif (anchor) {\n  // The three paragraphs right after the anchor are the ones to remove:\n  //   1) an empty spacer paragraph\n  //   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n  //   3) \"\u00a9 2020 . Contact: ... Creative Commons Attribution\"\n  const toDelete = [];\n  let current = anchor;\n  for (let i = 0; i < 3; i++) {\n    current = current.getNextOrNullObject();\n    current.load(\"text,isNullObject\");\n    await context.sync();\n    if (current.isNullObject) {\n      break;\n    }\n    toDelete.push(current);\n  }\n\n  // Sanity-check the expected text before deleting, then remove them in\n  // reverse order so earlier getNext() references stay valid.\n  const expected = [\n    \"\",\n    \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n    \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n  ];\n  for (let i = 0; i < toDelete.length; i++) {\n    if (toDelete[i].text !== expected[i]) {\n      throw new Error(\"Unexpected paragraph text while locating block to delete: \" + toDelete[i].text);\n    }\n  }\n\n  for (let i = toDelete.length - 1; i >= 0; i--) {\n    toDelete[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter / Salvar em pdf / Salvar em docx\" line,\n# the \"\u00a9 2020 ... Creative Commons Attribution\" footer line, and the blank\n# paragraph that separated them from the bibliography, while leaving the\n# rest of the document (including the FLEMMING bibliography entry and the\n# paragraph(s) that follow the footer) untouched.\n\n$d = $word.ActiveDocument\n\n# Locate the bibliography entry that immediately precedes the block we\n# need to drop.\n$anchor = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*FLEMMING, Diva M.; GON\u00c7ALVES, Mirian B.*Pearson Prentice Hall, 2009.*\") {\n        $anchor = $p\n        break\n    }\n}\n\nif ($anchor -ne $null) {\n    # The three paragraphs right after the anchor are the ones to remove:\n    #   1) an empty spacer paragraph\n    #   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n    #   3) \"\u00a9 2020 . Contact: ... Creative Commons Attribution\"\n    $expected = @(\n        \"\",\n        \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n        \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n    )\n\n    for ($i = 0; $i -lt $expected.Length; $i++) {\n        $target = $anchor.Next()\n        $actual = $target.Range.Text.TrimEnd()\n        if ($actual -ne $expected[$i]) {\n            throw \"Unexpected paragraph text while locating block to delete: \" + $actual\n        }\n        $target.Range.Delete()\n    }\n}\n"}
